# ---------------------------------------------------------------------------
# Dev V2 edit: add a new "Sheet2" worksheet (corner/block detection demo) in
# front of the existing "Sheet1", and extend "Sheet1" with a second copy of
# the same demo block further down the sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Insert a brand new worksheet named "Sheet2" before "Sheet1" ---------
$sheet1Ref = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 1, [System.Reflection.Missing]::Value)
$newSheet.Name = "Sheet2"
$newSheet.Move($sheet1Ref)

# Re-fetch fresh references by name now that the sheet collection has
# settled, since worksheet references obtained before a Move/Add can end
# up pointing at the wrong sheet once positions shift.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Use an already-bordered cell from "Sheet1" as the paint source for the
# bordered grid cells (this reuses the existing cell style instead of
# creating a new one).
$borderSource = $sheet1.Range("C5")

$borderSource.Copy()
$ws2.Range("D8:J14").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

# top / bottom header labels of the box
$ws2.Range("D8").Value = "F"
$ws2.Range("J8").Value = "G"
$ws2.Range("D14").Value = "H"
$ws2.Range("J14").Value = "J"

# numeric rows inside the box
for ($r = 9; $r -le 13; $r++) {
    for ($i = 0; $i -lt 5; $i++) {
        $col = [char](69 + $i)   # E,F,G,H,I
        $ws2.Range("$col$r").Value = $i + 1
    }
}

# "leaf" annotations around the box
$ws2.Range("K7").Value = "leaf"
$ws2.Range("L8").Value = "leaf 2"
$ws2.Range("L13").Value = "OMG Leaf"
$ws2.Range("L14").Value = "leaf 2"
$ws2.Range("D16").Value = "leaf 2"
$ws2.Range("C17").Value = "leaf 3"
$ws2.Range("C15").Value = "leaf 1"
$ws2.Range("K15").Value = "leaf"
$ws2.Range("M16").Value = "disconnected"

# --- 3. Add a second copy of the same block further down "Sheet1" ----------
$borderSource.Copy()
$sheet1.Range("E23:K30").PasteSpecial(-4122)
$sheet1.Application.CutCopyMode = $false

$sheet1.Range("E23").Value = "F"
$sheet1.Range("K23").Value = "G"
$sheet1.Range("E30").Value = "H"
$sheet1.Range("K30").Value = "J"

for ($r = 24; $r -le 29; $r++) {
    for ($i = 0; $i -lt 5; $i++) {
        $col = [char](70 + $i)   # F,G,H,I,J
        $sheet1.Range("$col$r").Value = $i + 1
    }
}

$sheet1.Range("L22").Value = "leaf"
$sheet1.Range("M23").Value = "leaf 2"
$sheet1.Range("M29").Value = "OMG Leaf"
$sheet1.Range("M30").Value = "leaf 2"
$sheet1.Range("E32").Value = "leaf 2"
$sheet1.Range("D33").Value = "leaf 3"
$sheet1.Range("D31").Value = "leaf 1"
$sheet1.Range("L31").Value = "leaf"
$sheet1.Range("N32").Value = "disconnected"

# --- 4. Restore the view / selection state ----------------------------------
$sheet1.Range("D21:N33").Select()
$sheet1.Activate()

$ws2.Range("J22").Select()
$ws2.Activate()
